# Swap the data that differ between row 4 and row 5 (Id, Antal, Enhet,
# Ost, Nord, Externid, Starttid, Sluttid) - the two observation records
# had their coordinates / counts mixed up and this corrects the mapping.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that actually differ between row 4 and row 5: A, I, J, Q, R, X, Z, AB
$cols = @("A", "I", "J", "Q", "R", "X", "Z", "AB")

foreach ($col in $cols) {
    $cell4 = $ws.Range("${col}4")
    $cell5 = $ws.Range("${col}5")

    $val4 = $cell4.Value2
    $val5 = $cell5.Value2

    $cell4.Value2 = $val5
    $cell5.Value2 = $val4
}
